$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1178
$ws1.Range("F3").Value = 1094
$ws1.Range("F4").Value = 1897
$ws1.Range("F5").Value = 594
$ws1.Range("F6").Value = 1241
$ws1.Range("F8").Value = 30
$ws1.Range("F9").Value = 134
$ws1.Range("F11").Value = 111
$ws1.Range("F12").Value = 98
$ws1.Range("F13").Value = 791
$ws1.Range("G13").Value = 68
$ws1.Range("F14").Value = 226
$ws1.Range("F15").Value = 121
$ws1.Range("F19").Value = 212
$ws1.Range("F20").Value = 690
$ws1.Range("F21").Value = 67
$ws1.Range("F23").Value = 182
$ws1.Range("F25").Value = 901
$ws1.Range("F26").Value = 347
$ws1.Range("F28").Value = 57
$ws1.Range("F29").Value = 300
$ws1.Range("F32").Value = 419

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 331
$ws2.Range("F7").Value = 262
$ws2.Range("F9").Value = 4

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 326

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 326
$ws4.Range("F3").Value = 1178
$ws4.Range("F4").Value = 1094
$ws4.Range("F5").Value = 1897
$ws4.Range("F6").Value = 594
$ws4.Range("F7").Value = 1241
$ws4.Range("F10").Value = 30
$ws4.Range("F11").Value = 134
$ws4.Range("F13").Value = 111
$ws4.Range("F14").Value = 98
$ws4.Range("F15").Value = 791
$ws4.Range("G15").Value = 68
$ws4.Range("F16").Value = 226
$ws4.Range("F17").Value = 121
$ws4.Range("F20").Value = 331
$ws4.Range("F25").Value = 262
$ws4.Range("F26").Value = 262
$ws4.Range("F27").Value = 212
$ws4.Range("F28").Value = 690
$ws4.Range("F29").Value = 67
$ws4.Range("F31").Value = 182
$ws4.Range("F33").Value = 901
$ws4.Range("F34").Value = 347
$ws4.Range("F36").Value = 4
$ws4.Range("F38").Value = 57
$ws4.Range("F39").Value = 300
$ws4.Range("F46").Value = 419
